# Finally was able to get single parameter sensitivity analysis to work.
#
# On the "mads" sheet, a couple of the LOG10(...) sensitivity formulas were
# replaced with their plain computed values (no longer live formulas), and
# the current cell selection moved from B11 to D9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mads")

# Replace formulas with static values (formula removed, value kept/updated).
$ws.Range("B2").Value = -5
$ws.Range("D2").Value = -3
$ws.Range("D4").Value = -3

# Update the active sheet/selection to D9 on the "mads" sheet.
$ws.Activate()
$ws.Range("D9").Select()
